# Update gh-pages output generated at 456a3b4
# - Refresh "想去人数" (F column) counters on several existing rows
# - Append a newly scraped event ("苏州·第三届华盟国漫次元嘉年华") as a new
#   last row on both the "展览" sheet and the "全部类型" sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions) — rows are 1-indexed, row 1 is the header
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value2  = 2097
$ws1.Range("F6").Value2  = 654
$ws1.Range("F8").Value2  = 2084
$ws1.Range("F9").Value2  = 10809
$ws1.Range("F15").Value2 = 9022
$ws1.Range("F16").Value2 = 1118
$ws1.Range("F17").Value2 = 731
$ws1.Range("F18").Value2 = 5298
$ws1.Range("F20").Value2 = 3360

# Append new row 21, matching the formatting of the previous last row (20)
$ws1.Range("A20").Copy()
$ws1.Range("A21").PasteSpecial(-4122)

$ws1.Range("A21").Value2 = 20
# Force column B to stay plain text (otherwise the "yyyy-mm-dd"-shaped
# string gets auto-parsed into a date serial number + date format, like
# every other date cell on this sheet, which is already text)
$ws1.Range("B21").NumberFormat = "@"
$ws1.Range("B21").Value2 = "2024-10-26"
$ws1.Range("B21").Style = "Normal"
$ws1.Range("C21").Value2 = "苏州·第三届华盟国漫次元嘉年华"
$ws1.Range("D21").Value2 = "清禾路886号 苏州聚橙尹山湖大剧院"
$ws1.Range("E21").Value2 = "2024.10.26 10:00-10.27 17:00"
$ws1.Range("F21").Value2 = 0
$ws1.Range("G21").Value2 = "不可售"
$ws1.Range("H21").Value2 = "https://show.bilibili.com/platform/detail.html?id=85767"
$ws1.Range("I21").Value2 = "//i1.hdslb.com/bfs/openplatform/202405/CqSYBZhQ1715846719965.jpeg"

# ---------------------------------------------------------------------
# Sheet "全部类型" (all categories) — same event set, different row offsets
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("全部类型")

$ws2.Range("F4").Value2  = 2097
$ws2.Range("F6").Value2  = 654
$ws2.Range("F9").Value2  = 2084
$ws2.Range("F12").Value2 = 10809
$ws2.Range("F18").Value2 = 9022
$ws2.Range("F19").Value2 = 1118
$ws2.Range("F20").Value2 = 731
$ws2.Range("F21").Value2 = 5298
$ws2.Range("F23").Value2 = 3360

# Append new row 24, matching the formatting of the previous last row (23)
$ws2.Range("A23").Copy()
$ws2.Range("A24").PasteSpecial(-4122)

$ws2.Range("A24").Value2 = 23
# Same text-coercion guard as sheet "展览" above
$ws2.Range("B24").NumberFormat = "@"
$ws2.Range("B24").Value2 = "2024-10-26"
$ws2.Range("B24").Style = "Normal"
$ws2.Range("C24").Value2 = "苏州·第三届华盟国漫次元嘉年华"
$ws2.Range("D24").Value2 = "清禾路886号 苏州聚橙尹山湖大剧院"
$ws2.Range("E24").Value2 = "2024.10.26 10:00-10.27 17:00"
$ws2.Range("F24").Value2 = 0
$ws2.Range("G24").Value2 = "不可售"
$ws2.Range("H24").Value2 = "https://show.bilibili.com/platform/detail.html?id=85767"
$ws2.Range("I24").Value2 = "//i1.hdslb.com/bfs/openplatform/202405/CqSYBZhQ1715846719965.jpeg"
